$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Typography")
$ws2 = $wb.Worksheets.Item("Translation")

# --- Typography sheet: add a new row (row 5) to the Typography table ---
# Pre-clear style on the whole new row range so freshly-created cells don't
# inherit the column's default style index (matches how the row already
# looks for untouched cells such as I5).
$ws1.Range("B5:J5").Style = "Normal"

$ws1.Range("B5").Value = "Typography_00"
$ws1.Range("C5").Value = "Asap-Regular.ttf"
$ws1.Range("D5").Value = 20
$ws1.Range("E5").Value = 4
$ws1.Range("F5").Value = "?"
$ws1.Range("G5").Value = "``~!@ #$%^&*()_+-={}|[]\<>?,./`"';:"
$ws1.Range("H5").Value = "0-9,a-z,A-Z"
# I5 already exists blank; J5 needs to be materialized as a blank cell too.
$ws1.Range("J5").Font.Bold = $false

# --- Translation sheet: update existing rows 4-6 and add new button/status-bar rows 7-12 ---

# Row 4: ResourceId1 / Default / Left / "Disco DCC LCC" / LTR
$ws2.Range("B4").Value = "ResourceId1"
$ws2.Range("D4").Value = "Left"
$ws2.Range("E4").Value = "Disco DCC LCC"

# Row 5: wildcardTextId / Typography_00 / Left / "*" / LTR
$ws2.Range("B5").Value = "wildcardTextId"
$ws2.Range("C5").Value = "Typography_00"
$ws2.Range("E5").Value = "*"

# Row 6: SingleUseId7 / Typography_00 / Center / "<value>" / LTR
$ws2.Range("B6").Value = "SingleUseId7"
$ws2.Range("C6").Value = "Typography_00"
$ws2.Range("D6").Value = "Center"
$ws2.Range("E6").Value = "<value>"

# Rows 7-12 are brand new - clear inherited column style first.
$ws2.Range("B7:F12").Style = "Normal"

# Row 7 (new ui button)
$ws2.Range("B7").Value = "SingleUseId8"
$ws2.Range("C7").Value = "Typography_00"
$ws2.Range("D7").Value = "Center"
$ws2.Range("E7").Value = "<value>"
$ws2.Range("F7").Value = "LTR"

# Row 8 (new ui button)
$ws2.Range("B8").Value = "SingleUseId9"
$ws2.Range("C8").Value = "Typography_00"
$ws2.Range("D8").Value = "Center"
$ws2.Range("E8").Value = "<value>"
$ws2.Range("F8").Value = "LTR"

# Row 9 (new ui button)
$ws2.Range("B9").Value = "SingleUseId10"
$ws2.Range("C9").Value = "Typography_00"
$ws2.Range("D9").Value = "Center"
$ws2.Range("E9").Value = "<value>"
$ws2.Range("F9").Value = "LTR"

# Row 10 (new ui button)
$ws2.Range("B10").Value = "SingleUseId11"
$ws2.Range("C10").Value = "Typography_00"
$ws2.Range("D10").Value = "Center"
$ws2.Range("E10").Value = "<value>"
$ws2.Range("F10").Value = "LTR"

# Row 11 (new ui button)
$ws2.Range("B11").Value = "SingleUseId12"
$ws2.Range("C11").Value = "Typography_00"
$ws2.Range("D11").Value = "Center"
$ws2.Range("E11").Value = "<value>"
$ws2.Range("F11").Value = "LTR"

# Row 12 (new ui button, link to preferences)
$ws2.Range("B12").Value = "SingleUseId13"
$ws2.Range("C12").Value = "Typography_00"
$ws2.Range("D12").Value = "Left"
$ws2.Range("E12").Value = "<value>"
$ws2.Range("F12").Value = "LTR"

Write-Host "Edits applied"
